# This script rearranges the columns of the single worksheet in the workbook.
# Old layout (A..M): AREA_CODE, AREA_NAME, BRANCH_CODE, BRANCH_NAME, RM_CODE, RM_NAME,
#                     BST_CODE, BST_NAME, SEG_ID, SEG_NAME, OS_TARGET_AMT, DISB_TARGET_AMT, INC_TARGET_AMT
# New layout (A..M): SEG_ID, SEG_NAME, AREA_CODE, AREA_NAME, BRANCH_CODE, BRANCH_NAME, RM_CODE, RM_NAME,
#                     BST_CODE, BST_NAME, OS_TARGET_AMT, DISB_TARGET_AMT, INC_TARGET_AMT
# i.e. the SEG_ID/SEG_NAME columns (old I:J) are moved in front of AREA_CODE (new A:B),
# and everything that used to be in A:H shifts right into C:J. K:M (the target amount columns) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- Step 1: stage a full copy (values + formats) of the A1:J4 block far out of the way (P1:Y4) ---
$ws.Range("A1:J4").Copy()
$ws.Range("P1").PasteSpecial($xlPasteValues)
$ws.Range("A1:J4").Copy()
$ws.Range("P1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Staged columns: P=old A, Q=old B, R=old C, S=old D, T=old E, U=old F, V=old G, W=old H, X=old I, Y=old J

# --- Step 2: rebuild A1:J4 from the staged copy, in the new column order ---
# new destination column -> staged source column
$pairs = @(
  @("A", "X"),  # SEG_ID        <- old I
  @("B", "Y"),  # SEG_NAME      <- old J
  @("C", "P"),  # AREA_CODE     <- old A
  @("D", "Q"),  # AREA_NAME     <- old B
  @("E", "R"),  # BRANCH_CODE   <- old C
  @("F", "S"),  # BRANCH_NAME   <- old D
  @("G", "T"),  # RM_CODE       <- old E
  @("H", "U"),  # RM_NAME       <- old F
  @("I", "V"),  # BST_CODE      <- old G
  @("J", "W")   # BST_NAME      <- old H
)

foreach ($p in $pairs) {
  $dst = $p[0]
  $src = $p[1]
  $dstRange = "$dst" + "1:" + "$dst" + "4"
  $srcRange = "$src" + "1:" + "$src" + "4"

  # clear destination contents first so that blank source cells really end up blank
  $ws.Range($dstRange).ClearContents()

  $ws.Range($srcRange).Copy()
  $ws.Range($dst + "1").PasteSpecial($xlPasteValues)

  $ws.Range($srcRange).Copy()
  $ws.Range($dst + "1").PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = 0

# --- Step 3: clear the staging area completely (contents + formatting) ---
$ws.Range("P1:Y4").Clear()

# --- Step 4: fix up the sheet view (selection moved, topLeftCell scroll reset) ---
$window = $excel.ActiveWindow
$window.ScrollColumn = 1
[void]$ws.Range("C14").Select()
